# Updated to 11_12 measurements
$wb = $excel.ActiveWorkbook

# --- 1. Fill in the previously-blank dilution-rate readings on sheet "11_5" ---
$ws5 = $wb.Worksheets.Item("11_5")

$ws5.Range("A4").Formula = "=24+7+3/4"
$ws5.Range("B4").Value = 3.55

$ws5.Range("A5").Formula = "=24+40/60"
$ws5.Range("B5").Value = 3.125

$ws5.Range("A6").Formula = "=48+20+0.25"
$ws5.Range("B6").Value = 8.35

# --- 2. Duplicate the "11_5" sheet to create the new "11_12" sheet ---
# (Copy places the clone right after $ws5 and activates it, mirroring Excel's UI behavior)
$ws5.Copy($null, $ws5)
$wsNew = $wb.ActiveSheet
$wsNew.Name = "11_12"

# --- 3. Update the new sheet's measurements for 11_12 ---
$wsNew.Range("B1").Value = 0.83

$wsNew.Range("A12").Value = 12.8
$wsNew.Range("F12").Value = 10464

$wsNew.Range("A13").Value = 12.81
$wsNew.Range("F13").Value = 10676

$wsNew.Range("A14").Value = 12.64
$wsNew.Range("F14").Value = 10245

$wsNew.Range("A15").Value = 12.78

$wsNew.Range("A16").Value = 12.78

# --- 4. Fix up the view state on both sheets ---
$ws5.Activate()
$ws5.Range("A12").Select()

$wsNew.Activate()
$wsNew.Range("F15").Select()

$wb.Save()
